# Apply weekly update of fruit/vegetable (Haba) prices.
# The data rows (2-14) are effectively reshuffled/updated: date (D), volume (J),
# min/max/avg price (K/L/M), origin (O) and $/Kg price (P) columns change values
# for each row as described in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value  = 45113
$ws.Cells.Item(2, 10).Value = 8
$ws.Cells.Item(2, 11).Value = 17000
$ws.Cells.Item(2, 12).Value = 17000
$ws.Cells.Item(2, 13).Value = 17000
$ws.Cells.Item(2, 16).Value = 680

# Row 3
$ws.Cells.Item(3, 4).Value  = 44425
$ws.Cells.Item(3, 10).Value = 25
$ws.Cells.Item(3, 11).Value = 14000
$ws.Cells.Item(3, 12).Value = 14000
$ws.Cells.Item(3, 13).Value = 14000
$ws.Cells.Item(3, 16).Value = 560

# Row 4
$ws.Cells.Item(4, 4).Value  = 44446
$ws.Cells.Item(4, 10).Value = 15
$ws.Cells.Item(4, 11).Value = 13000
$ws.Cells.Item(4, 12).Value = 13000
$ws.Cells.Item(4, 13).Value = 13000
$ws.Cells.Item(4, 16).Value = 520

# Row 5
$ws.Cells.Item(5, 4).Value  = 44453
$ws.Cells.Item(5, 10).Value = 55
$ws.Cells.Item(5, 11).Value = 14000
$ws.Cells.Item(5, 12).Value = 15000
$ws.Cells.Item(5, 13).Value = 14455
$ws.Cells.Item(5, 16).Value = 578

# Row 6
$ws.Cells.Item(6, 4).Value  = 44376
$ws.Cells.Item(6, 11).Value = 12000
$ws.Cells.Item(6, 12).Value = 12000
$ws.Cells.Item(6, 13).Value = 12000
$ws.Cells.Item(6, 16).Value = 480

# Row 7
$ws.Cells.Item(7, 4).Value  = 44421
$ws.Cells.Item(7, 10).Value = 20
$ws.Cells.Item(7, 11).Value = 15000
$ws.Cells.Item(7, 12).Value = 15000
$ws.Cells.Item(7, 13).Value = 15000
$ws.Cells.Item(7, 16).Value = 600

# Row 8
$ws.Cells.Item(8, 15).Value = "Provincia de Limarí"

# Row 9
$ws.Cells.Item(9, 4).Value  = 44435
$ws.Cells.Item(9, 10).Value = 15
$ws.Cells.Item(9, 11).Value = 14000
$ws.Cells.Item(9, 12).Value = 14000
$ws.Cells.Item(9, 13).Value = 14000
$ws.Cells.Item(9, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(9, 16).Value = 560

# Row 10
$ws.Cells.Item(10, 4).Value  = 44418
$ws.Cells.Item(10, 10).Value = 12
$ws.Cells.Item(10, 11).Value = 15000
$ws.Cells.Item(10, 12).Value = 15000
$ws.Cells.Item(10, 13).Value = 15000
$ws.Cells.Item(10, 16).Value = 600

# Row 11
$ws.Cells.Item(11, 4).Value  = 44340
$ws.Cells.Item(11, 10).Value = 25
$ws.Cells.Item(11, 11).Value = 15000
$ws.Cells.Item(11, 12).Value = 15000
$ws.Cells.Item(11, 13).Value = 15000
$ws.Cells.Item(11, 16).Value = 600

# Row 12
$ws.Cells.Item(12, 4).Value  = 44432
$ws.Cells.Item(12, 10).Value = 15
$ws.Cells.Item(12, 11).Value = 14000
$ws.Cells.Item(12, 12).Value = 14000
$ws.Cells.Item(12, 13).Value = 14000
$ws.Cells.Item(12, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(12, 16).Value = 560

# Row 13
$ws.Cells.Item(13, 4).Value  = 44449
$ws.Cells.Item(13, 10).Value = 30
$ws.Cells.Item(13, 11).Value = 16000
$ws.Cells.Item(13, 12).Value = 16000
$ws.Cells.Item(13, 13).Value = 16000
$ws.Cells.Item(13, 16).Value = 640

# Row 14
$ws.Cells.Item(14, 4).Value  = 44467
$ws.Cells.Item(14, 10).Value = 35
$ws.Cells.Item(14, 11).Value = 12000
$ws.Cells.Item(14, 12).Value = 12000
$ws.Cells.Item(14, 13).Value = 12000
$ws.Cells.Item(14, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(14, 16).Value = 480
